# Updates the cryptos price list (sheet1) cell values to match the
# latest scrape: D = Price text, E = Volume(1h) text; some coin rows
# (49/50) swapped position (Arweave <-> ONDO) with new data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.854.19"
$ws.Range("E2").Value = "  -3.93%  "
# Row 3
$ws.Range("D3").Value = "3.342.94"
$ws.Range("E3").Value = "  -0.89%  "
# Row 4
$ws.Range("E4").Value = "  -0.01%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.36"
$ws.Range("E5").Value = "  -3.32%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.66"
$ws.Range("E6").Value = "  -5.17%  "
# Row 7
$ws.Range("E7").Value = "  +0.02%  "
# Row 8
$ws.Range("E8").Value = "  -1.48%  "
# Row 9
$ws.Range("E9").Value = "  -3.68%  "
# Row 10
$ws.Range("E10").Value = "  -1.84%  "
# Row 11
$ws.Range("E11").Value = "  -4.37%  "
# Row 12
$ws.Range("D12").Value = "3.926.29"
$ws.Range("E12").Value = "  -0.87%  "
# Row 13
$ws.Range("E13").Value = "  -1.66%  "
# Row 14
$ws.Range("E14").Value = "  -5.16%  "
# Row 15
$ws.Range("D15").Value = "66.888.02"
# Row 16
$ws.Range("E16").Value = "  -2.34%  "
# Row 17
$ws.Range("D17").Value = "3.335.41"
$ws.Range("E17").Value = "  -0.32%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "436.16"
$ws.Range("E18").Value = "  -3.46%  "
# Row 19
$ws.Range("E19").Value = "  -1.45%  "
# Row 20
$ws.Range("E20").Value = "  -2.69%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.63"
$ws.Range("E21").Value = "  -2.34%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.88"
$ws.Range("E22").Value = "  -1.13%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.02%  "
# Row 24
$ws.Range("E24").Value = "  +0.18%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000118"
$ws.Range("E25").Value = "  -3.95%  "
# Row 26
$ws.Range("E26").Value = "  -0.32%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.04"
$ws.Range("E27").Value = "  -5.26%  "
# Row 28
$ws.Range("E28").Value = "  +0.09%  "
# Row 29
$ws.Range("E29").Value = "  -1.73%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.87"
$ws.Range("E30").Value = "  -1.93%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.32"
$ws.Range("E31").Value = "  -5.91%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.998"
$ws.Range("E32").Value = "  -0.02%  "
# Row 33
$ws.Range("E33").Value = "  -3.78%  "
# Row 34
$ws.Range("E34").Value = "  -3.07%  "
# Row 35
$ws.Range("E35").Value = "  -2.39%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.12"
$ws.Range("E36").Value = "  -2.36%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "27.87"
$ws.Range("E37").Value = "  +2.32%  "
# Row 38
$ws.Range("E38").Value = "  -4.77%  "
# Row 39
$ws.Range("D39").Value = "2.840.13"
$ws.Range("E39").Value = "  +3.72%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.792"
$ws.Range("E40").Value = "  -3.31%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.45"
$ws.Range("E41").Value = "  -3.67%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.29"
$ws.Range("E42").Value = "  -4.15%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0673"
$ws.Range("E43").Value = "  -2.69%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.18"
$ws.Range("E44").Value = "  -1.24%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "24.65"
$ws.Range("E45").Value = "  -4.17%  "
# Row 46
$ws.Range("E46").Value = "  -7.02%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "324.09"
$ws.Range("E47").Value = "  -5.54%  "
# Row 48
$ws.Range("E48").Value = "  -4.48%  "
# Row 49
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.987"
$ws.Range("E49").Value = "  -3.62%  "
# Row 50
$ws.Range("B50").Value = "Arweave"
$ws.Range("C50").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "31.31"
$ws.Range("E50").Value = "  -5.20%  "
# Row 51
$ws.Range("E51").Value = "  -2.77%  "
